$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 137 ---
$ws.Range("A136").Copy()
$ws.Range("A137").PasteSpecial(-4122)
$ws.Range("A137").Value = 44867
$ws.Range("C137").Value = "test"
$ws.Range("B137").Value = "MRN:  JH16121935"
$ws.Range("D137").Value = 0
$ws.Range("E137").Value = 20
$ws.Range("F137").Value = 32.901000000000003
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 1
$ws.Range("J137").Value = 1
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 1
$ws.Range("M137").Value = 0
$ws.Range("N137").Value = 0
$ws.Range("O137").Value = 0
$ws.Range("P137").Value = 0
$ws.Range("Q137").Value = 0
$ws.Range("R137").Value = 0
$ws.Range("S137").Value = 0
$ws.Range("T137").Value = 0
$ws.Range("U137").Value = 0
$ws.Range("V137").Value = 0
$ws.Range("W137").Value = 0
$ws.Range("X137").Value = 0
$ws.Range("Y137").Value = 0
$ws.Range("Z137").Value = 0
$ws.Range("AA137").Value = 0
$ws.Range("AB137").Value = 0

# --- Row 138 ---
$ws.Range("A136").Copy()
$ws.Range("A138").PasteSpecial(-4122)
$ws.Range("A138").Value = 44867
$ws.Range("B138").Value = "MRN:  JH16121937"
$ws.Range("C138").Value = "Khashab"
$ws.Range("D138").Value = 0
$ws.Range("E138").Value = 60
$ws.Range("F138").Value = 15.31
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 1
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 0
$ws.Range("N138").Value = 1
$ws.Range("O138").Value = 1
$ws.Range("P138").Value = 1
$ws.Range("Q138").Value = 0
$ws.Range("R138").Value = 1
$ws.Range("S138").Value = 0
$ws.Range("T138").Value = 0
$ws.Range("U138").Value = 1
$ws.Range("V138").Value = 0
$ws.Range("W138").Value = 0
$ws.Range("X138").Value = 0
$ws.Range("Y138").Value = 1
$ws.Range("Z138").Value = 0
$ws.Range("AA138").Value = 0
$ws.Range("AB138").Value = 0

# --- Restore view state: scrolled down near the bottom of the data, with
#     the next empty row (D145) selected, matching the author's final
#     on-screen position after entering the new rows. ---
$excel.ActiveWindow.ScrollRow = 110
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D145").Select() | Out-Null
